$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in virus_genus value (column C, row 2): "Icthamaparvovirus" -> "Ichthamaparvovirus"
$ws.Range("C2").Value = "Ichthamaparvovirus"

# Reflect the updated cell selection left by the author's edit
$ws.Range("C2").Select()
